$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates (shared strings for this row get registered first: indices 13-15)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = ";0"
$ws.Range("G4").Value = ";21"
$ws.Range("H4").Value = ";-1000"

# Row 3 updates (shared strings for this row get registered next: indices 16-18)
$ws.Range("B3").Value = 400
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = ";1"
$ws.Range("G3").Value = ";20"
$ws.Range("H3").Value = ";+100"

# Move selection/active cell to D4
$ws.Range("D4").Select()
